# Auto-generated script applying the Ridill_Profits market-data refresh.
# Updates columns H-N (currentAveragePrice .. LeveProfitHQ) for specific rows
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets to reflect refreshed market data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 137
$ws.Range("H137").Value = 19010990
$ws.Range("I137").Value = 28410108
$ws.Range("J137").Value = 11625968
$ws.Range("K137").Value = 85230324
$ws.Range("L137").Value = 34877904
$ws.Range("M137").Value = -85227774
$ws.Range("N137").Value = -34883004

# Row 141
$ws.Range("H141").Value = 3573.158
$ws.Range("I141").Value = 2023.3334
$ws.Range("J141").Value = 6230
$ws.Range("K141").Value = 6070.0002
$ws.Range("L141").Value = 18690
$ws.Range("M141").Value = -890.0002000000004
$ws.Range("N141").Value = -29050


$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 1459617.9
$ws.Range("I32").Value = 1629373.5
$ws.Range("J32").Value = 7264.4443
$ws.Range("K32").Value = 1629373.5
$ws.Range("L32").Value = 7264.4443
$ws.Range("M32").Value = -1629086.5
$ws.Range("N32").Value = -7838.4443

# Row 61
$ws.Range("H61").Value = 6082496.5
$ws.Range("I61").Value = 5809414.5
$ws.Range("J61").Value = 6537633
$ws.Range("K61").Value = 5809414.5
$ws.Range("L61").Value = 6537633
$ws.Range("M61").Value = -5809202.5
$ws.Range("N61").Value = -6538057

# Row 136
$ws.Range("H136").Value = 6082496.5
$ws.Range("I136").Value = 5809414.5
$ws.Range("J136").Value = 6537633
$ws.Range("K136").Value = 17428243.5
$ws.Range("L136").Value = 19612899
$ws.Range("M136").Value = -17425693.5
$ws.Range("N136").Value = -19617999


$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 13236341
$ws.Range("I134").Value = 16130074
$ws.Range("J134").Value = 3269036
$ws.Range("K134").Value = 48390222
$ws.Range("L134").Value = 9807108
$ws.Range("M134").Value = -48387687
$ws.Range("N134").Value = -9812178


$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 897.2727
$ws.Range("I16").Value = 826.7143
$ws.Range("J16").Value = 1020.75
$ws.Range("K16").Value = 826.7143
$ws.Range("L16").Value = 1020.75
$ws.Range("M16").Value = -539.7143
$ws.Range("N16").Value = -1594.75

# Row 31
$ws.Range("H31").Value = 1778.95
$ws.Range("I31").Value = 1126.1333
$ws.Range("J31").Value = 2313.0728
$ws.Range("K31").Value = 1126.1333
$ws.Range("L31").Value = 2313.0728
$ws.Range("M31").Value = -831.1333
$ws.Range("N31").Value = -2903.0728

# Row 34
$ws.Range("H34").Value = 1778.95
$ws.Range("I34").Value = 1126.1333
$ws.Range("J34").Value = 2313.0728
$ws.Range("K34").Value = 1126.1333
$ws.Range("L34").Value = 2313.0728
$ws.Range("M34").Value = -924.1333
$ws.Range("N34").Value = -2717.0728

# Row 113
$ws.Range("H113").Value = 897.2727
$ws.Range("I113").Value = 826.7143
$ws.Range("J113").Value = 1020.75
$ws.Range("K113").Value = 826.7143
$ws.Range("L113").Value = 1020.75
$ws.Range("M113").Value = 1343.2857
$ws.Range("N113").Value = -5360.75

# Row 122
$ws.Range("H122").Value = 4133.8276
$ws.Range("I122").Value = 5484.15
$ws.Range("J122").Value = 1133.1111
$ws.Range("K122").Value = 16452.45
$ws.Range("L122").Value = 3399.3333
$ws.Range("M122").Value = -14002.45
$ws.Range("N122").Value = -8299.3333


$ws = $wb.Worksheets.Item("CUL")
# Row 44
$ws.Range("H44").Value = 953.4
$ws.Range("I44").Value = 826.6667
$ws.Range("J44").Value = 1007.7143
$ws.Range("K44").Value = 2480.0001
$ws.Range("L44").Value = 3023.1429
$ws.Range("M44").Value = -2082.0001
$ws.Range("N44").Value = -3819.1429

# Row 46
$ws.Range("H46").Value = 1786.1034
$ws.Range("I46").Value = 1244.1666
$ws.Range("J46").Value = 1927.4783
$ws.Range("K46").Value = 3732.4998
$ws.Range("L46").Value = 5782.4349
$ws.Range("M46").Value = -3641.4998
$ws.Range("N46").Value = -5964.4349

# Row 51
$ws.Range("H51").Value = 2183.3333
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 2183.3333
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 6549.999899999999
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -7469.999899999999

# Row 57
$ws.Range("H57").Value = 3475
$ws.Range("I57").Value = 1950
$ws.Range("J57").Value = 5000
$ws.Range("K57").Value = 5850
$ws.Range("L57").Value = 15000
$ws.Range("M57").Value = -5291
$ws.Range("N57").Value = -16118

# Row 58
$ws.Range("H58").Value = 3100
$ws.Range("J58").Value = 3625
$ws.Range("L58").Value = 10875
$ws.Range("N58").Value = -11131

# Row 68
$ws.Range("H68").Value = 1997.6263
$ws.Range("I68").Value = 555.2646999999999
$ws.Range("J68").Value = 2857.9824
$ws.Range("K68").Value = 1665.7941
$ws.Range("L68").Value = 8573.947199999999
$ws.Range("M68").Value = -854.7940999999998
$ws.Range("N68").Value = -10195.9472

# Row 71
$ws.Range("H71").Value = 1997.6263
$ws.Range("I71").Value = 555.2646999999999
$ws.Range("J71").Value = 2857.9824
$ws.Range("K71").Value = 4997.382299999999
$ws.Range("L71").Value = 25721.8416
$ws.Range("M71").Value = -941.3822999999993
$ws.Range("N71").Value = -33833.8416

# Row 113
$ws.Range("H113").Value = 2386.7917
$ws.Range("I113").Value = 3298.0833
$ws.Range("J113").Value = 1475.5
$ws.Range("K113").Value = 9894.249899999999
$ws.Range("L113").Value = 4426.5
$ws.Range("M113").Value = -7724.249899999999
$ws.Range("N113").Value = -8766.5


$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 19325590
$ws.Range("I132").Value = 15334398
$ws.Range("J132").Value = 25977580
$ws.Range("K132").Value = 46003194
$ws.Range("L132").Value = 77932740
$ws.Range("M132").Value = -46000664
$ws.Range("N132").Value = -77937800


$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 1443.3478
$ws.Range("I7").Value = 1322.35
$ws.Range("J7").Value = 2250
$ws.Range("K7").Value = 1322.35
$ws.Range("L7").Value = 2250
$ws.Range("M7").Value = -1210.35
$ws.Range("N7").Value = -2474

# Row 22
$ws.Range("H22").Value = 38471584
$ws.Range("I22").Value = 4240
$ws.Range("J22").Value = 55568180
$ws.Range("K22").Value = 4240
$ws.Range("L22").Value = 55568180
$ws.Range("M22").Value = -3945
$ws.Range("N22").Value = -55568770

# Row 27
$ws.Range("H27").Value = 38471584
$ws.Range("I27").Value = 4240
$ws.Range("J27").Value = 55568180
$ws.Range("K27").Value = 4240
$ws.Range("L27").Value = 55568180
$ws.Range("M27").Value = -4133
$ws.Range("N27").Value = -55568394

# Row 40
$ws.Range("H40").Value = 13890353
$ws.Range("I40").Value = 15874332
$ws.Range("K40").Value = 15874332
$ws.Range("M40").Value = -15874196

# Row 46
$ws.Range("H46").Value = 6160.8
$ws.Range("I46").Value = 10000
$ws.Range("J46").Value = 402
$ws.Range("K46").Value = 10000
$ws.Range("L46").Value = 402
$ws.Range("M46").Value = -9812
$ws.Range("N46").Value = -778

# Row 68
$ws.Range("H68").Value = 2134.7585
$ws.Range("I68").Value = 1817.1052
$ws.Range("J68").Value = 2738.3
$ws.Range("K68").Value = 1817.1052
$ws.Range("L68").Value = 2738.3
$ws.Range("M68").Value = -1068.1052
$ws.Range("N68").Value = -4236.3

# Row 71
$ws.Range("H71").Value = 2134.7585
$ws.Range("I71").Value = 1817.1052
$ws.Range("J71").Value = 2738.3
$ws.Range("K71").Value = 9085.526
$ws.Range("L71").Value = 13691.5
$ws.Range("M71").Value = -5341.526
$ws.Range("N71").Value = -21179.5

# Row 100
$ws.Range("H100").Value = 1707.5294
$ws.Range("I100").Value = 1466.4445
$ws.Range("J100").Value = 1978.75
$ws.Range("K100").Value = 1466.4445
$ws.Range("L100").Value = 1978.75
$ws.Range("M100").Value = -925.4445000000001
$ws.Range("N100").Value = -3060.75

# Row 122
$ws.Range("H122").Value = 26184306
$ws.Range("I122").Value = 18086966
$ws.Range("J122").Value = 66671000
$ws.Range("K122").Value = 54260898
$ws.Range("L122").Value = 200013000
$ws.Range("M122").Value = -54258448
$ws.Range("N122").Value = -200017900

# Row 126
$ws.Range("H126").Value = 1443.3478
$ws.Range("I126").Value = 1322.35
$ws.Range("J126").Value = 2250
$ws.Range("K126").Value = 3967.05
$ws.Range("L126").Value = 6750
$ws.Range("M126").Value = -1497.05
$ws.Range("N126").Value = -11690


$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 382.6154
$ws.Range("I113").Value = 249.83333
$ws.Range("J113").Value = 496.42856
$ws.Range("K113").Value = 749.49999
$ws.Range("L113").Value = 1489.28568
$ws.Range("M113").Value = 1420.50001
$ws.Range("N113").Value = -5829.28568

# Row 122
$ws.Range("H122").Value = 2416386.8
$ws.Range("I122").Value = 4116015.2
$ws.Range("J122").Value = 1125.2632
$ws.Range("K122").Value = 12348045.6
$ws.Range("L122").Value = 3375.7896
$ws.Range("M122").Value = -12345595.6
$ws.Range("N122").Value = -8275.7896
